$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.216.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "'2.228.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.58%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'318.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("D6").Value = "'98.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.85%  "
$ws.Range("E7").Value = "  -7.31%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  -7.51%  "
$ws.Range("D10").Value = "'37.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.47%  "
$ws.Range("D11").Value = "'54.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("D12").Value = "'0.0830"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.68%  "
$ws.Range("D13").Value = "'7.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.88%  "
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "'2.564.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.83%  "
$ws.Range("D16").Value = "'0.863"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -10.16%  "
$ws.Range("D17").Value = "'14.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.13%  "
$ws.Range("D18").Value = "'2.223.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.19%  "
$ws.Range("D19").Value = "'43.103.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").Value = "'13.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.41%  "
$ws.Range("E21").Value = "  -8.67%  "
$ws.Range("E22").Value = "  -8.05%  "
$ws.Range("D23").Value = "'3.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.74%  "
$ws.Range("D24").Value = "'65.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.98%  "
$ws.Range("D25").Value = "'237.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.98%  "
$ws.Range("D26").Value = "'2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").Value = "'10.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.54%  "
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("D31").Value = "'6.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.39%  "
$ws.Range("D32").Value = "'36.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "'20.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.05%  "
$ws.Range("D34").Value = "'0.0866"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.29%  "
$ws.Range("D35").Value = "'157.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.86%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.92%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'3.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.71%  "
$ws.Range("E38").Value = "  -7.18%  "
$ws.Range("D39").Value = "'1.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").Value = "'4.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.96%  "
$ws.Range("E41").Value = "  -9.15%  "
$ws.Range("D42").Value = "'3.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.50%  "
$ws.Range("D43").Value = "'0.0319"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.27%  "
$ws.Range("D44").Value = "'14.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.17%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'1.756.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.57%  "
$ws.Range("D47").Value = "'0.203"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.33%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'83.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.40%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("D50").Value = "'5.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -11.94%  "
$ws.Range("D51").Value = "'73.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.50%  "
